# Implemented r-skript for mca baeume
# The "barrierefreiheit" criterion row (row 10, part of the
# "multifunktionale_nutzungsqualitaet" group) is removed from the
# "Gewichtung" sheet. The remaining three criteria in that group
# are re-weighted evenly (0.333333333 instead of 0.25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gewichtung")

# Remove the "barrierefreiheit" criterion row entirely (old row 10);
# this shifts the following "kreislauffaehigkeit" rows up by one and
# shrinks the table from 13 to 12 data rows.
$ws.Rows("10").Delete()

# Re-balance the within-group weights for the remaining criteria of the
# "multifunktionale_nutzungsqualitaet" group (now rows 7-9) from 0.25 to
# an even three-way split.
$ws.Range("F7").Value = 0.333333333
$ws.Range("F8").Value = 0.333333333
$ws.Range("F9").Value = 0.333333333

# Update the active selection on the sheet to match the author's final
# cursor position.
$ws.Range("F16").Select()
